$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update statut (A) values to "4"/"1" (must stay text, not numeric) ---
$ws.Range("A2:A4").NumberFormat = "@"

$ws.Range("A2").Value = "4"
$ws.Range("A3").Value = "1"
$ws.Range("A4").Value = "1"

$ws.Range("A2:A4").ClearFormats()

# --- Update statut_name (C) values ---
$ws.Range("C2").Value = "4: pas de résultats postés ni publiés"
$ws.Range("C3").Value = "1: résultats postés ou publiés dans les 12 mois"
$ws.Range("C4").Value = "1: résultats postés ou publiés dans les 12 mois"

# --- Delete obsolete columns (right to left so letters don't shift) ---
$ws.Range("L:L").Delete()   # results
$ws.Range("K:K").Delete()   # results_3y
$ws.Range("J:J").Delete()   # results_1y
$ws.Range("B:B").Delete()   # statut_label
